# Fruta / hortaliza, semanal
# Insert a new weekly record as row 262 (pushing the existing rows 262-308
# down to 263-309) in the Berenjena - Terminal La Palmera de La Serena sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 262..308 down to 263..309, leaving a blank row 262 to fill in.
$ws.Rows(262).Insert()

# Populate the newly inserted row 262 with the new weekly price record.
$ws.Cells.Item(262, 1).Value2  = 8
$ws.Cells.Item(262, 2).Value2  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(262, 3).Value2  = "Coquimbo"
$ws.Cells.Item(262, 4).Value2  = 45209
$ws.Cells.Item(262, 5).Value2  = 4
$ws.Cells.Item(262, 6).Value2  = 100112001
$ws.Cells.Item(262, 7).Value2  = "Berenjena"
$ws.Cells.Item(262, 8).Value2  = "Sin especificar"
$ws.Cells.Item(262, 9).Value2  = "Primera"
$ws.Cells.Item(262, 10).Value2 = 440
$ws.Cells.Item(262, 11).Value2 = 9500
$ws.Cells.Item(262, 12).Value2 = 10000
$ws.Cells.Item(262, 13).Value2 = 9750
$ws.Cells.Item(262, 14).Value2 = "$/caja 50 unidades"
$ws.Cells.Item(262, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(262, 16).Value2 = 195
$ws.Cells.Item(262, 17).Value2 = 50
$ws.Cells.Item(262, 18).Value2 = "Hortaliza"
